# Auto-generated edit script
# Applies cell-level numeric updates to match the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (hunk 0)
$ws.Range("H11").Value = 1611.2858
$ws.Range("I11").Value = 1611.2858
$ws.Range("K11").Value = 1611.2858
$ws.Range("M11").Value = -1471.2858
# Row 40 (hunk 1)
$ws.Range("H40").Value = 3800.611
$ws.Range("I40").Value = 3082.1
$ws.Range("K40").Value = 3082.1
$ws.Range("M40").Value = -2907.1
# Row 64 (hunk 2)
$ws.Range("H64").Value = 6510.769
$ws.Range("I64").Value = 4900
$ws.Range("J64").Value = 6803.636
$ws.Range("K64").Value = 4900
$ws.Range("L64").Value = 6803.636
$ws.Range("M64").Value = -4652
$ws.Range("N64").Value = -7299.636
# Row 67 (hunk 3)
$ws.Range("H67").Value = 6510.769
$ws.Range("I67").Value = 4900
$ws.Range("J67").Value = 6803.636
$ws.Range("K67").Value = 4900
$ws.Range("L67").Value = 6803.636
$ws.Range("M67").Value = -4042
$ws.Range("N67").Value = -8519.636
# Row 97 (hunk 4)
$ws.Range("H97").Value = 2555
$ws.Range("J97").Value = 2555
$ws.Range("L97").Value = 7665
$ws.Range("N97").Value = -8657
# Row 112 (hunk 5)
$ws.Range("H112").Value = 2089.111
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").Value = ""
# Row 129 (hunk 6)
$ws.Range("H129").Value = 2019.5
$ws.Range("I129").Value = 1404.25
$ws.Range("J129").Value = 3250
$ws.Range("K129").Value = 4212.75
$ws.Range("L129").Value = 9750
$ws.Range("M129").Value = 787.25
$ws.Range("N129").Value = -19750
# Row 137 (hunk 7)
$ws.Range("H137").Value = 2700.5134
$ws.Range("I137").Value = 3063.423
$ws.Range("J137").Value = 1842.7273
$ws.Range("K137").Value = 9190.269
$ws.Range("L137").Value = 5528.1819
$ws.Range("M137").Value = -6640.269
$ws.Range("N137").Value = -10628.1819
# Row 138 (hunk 8)
$ws.Range("H138").Value = 1978.24
$ws.Range("I138").Value = 1031.7222
$ws.Range("J138").Value = 2186.0122
$ws.Range("K138").Value = 3095.1666
$ws.Range("L138").Value = 6558.0366
$ws.Range("M138").Value = 2044.8334
$ws.Range("N138").Value = -16838.0366

$ws = $wb.Worksheets.Item("BSM")
# Row 88 (hunk 9)
$ws.Range("H88").Value = 29996
$ws.Range("I88").Value = 29990
$ws.Range("J88").Value = 29999
$ws.Range("K88").Value = 29990
$ws.Range("L88").Value = 29999
$ws.Range("M88").Value = -29584
$ws.Range("N88").Value = -30811
# Row 91 (hunk 10)
$ws.Range("H91").Value = 29996
$ws.Range("I91").Value = 29990
$ws.Range("J91").Value = 29999
$ws.Range("K91").Value = 29990
$ws.Range("L91").Value = 29999
$ws.Range("M91").Value = -28586
$ws.Range("N91").Value = -32807
# Row 134 (hunk 11)
$ws.Range("H134").Value = 716949.6
$ws.Range("I134").Value = 2867.7693
$ws.Range("K134").Value = 8603.3079
$ws.Range("M134").Value = -6068.3079

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (hunk 12)
$ws.Range("H31").Value = 654599
$ws.Range("I31").Value = 9013
$ws.Range("J31").Value = 1796789.8
$ws.Range("K31").Value = 9013
$ws.Range("L31").Value = 1796789.8
$ws.Range("M31").Value = -8718
$ws.Range("N31").Value = -1797379.8
# Row 34 (hunk 13)
$ws.Range("H34").Value = 654599
$ws.Range("I34").Value = 9013
$ws.Range("J34").Value = 1796789.8
$ws.Range("K34").Value = 9013
$ws.Range("L34").Value = 1796789.8
$ws.Range("M34").Value = -8811
$ws.Range("N34").Value = -1797193.8
# Row 43 (hunk 14)
$ws.Range("H43").Value = 44137.668
$ws.Range("J43").Value = 44137.668
$ws.Range("L43").Value = 44137.668
$ws.Range("N43").Value = -44505.668
# Row 101 (hunk 15)
$ws.Range("H101").Value = 44137.668
$ws.Range("J101").Value = 44137.668
$ws.Range("L101").Value = 44137.668
$ws.Range("N101").Value = -50627.668
# Row 107 (hunk 16)
$ws.Range("H107").Value = 2249.842
$ws.Range("I107").Value = 849.1429000000001
$ws.Range("J107").Value = 3066.9167
$ws.Range("K107").Value = 849.1429000000001
$ws.Range("L107").Value = 3066.9167
$ws.Range("M107").Value = 1070.8571
$ws.Range("N107").Value = -6906.9167

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (hunk 17)
$ws.Range("H5").Value = 1833.6316
$ws.Range("I5").Value = 1889.9375
$ws.Range("K5").Value = 5669.8125
$ws.Range("M5").Value = -5557.8125
# Row 18 (hunk 18)
$ws.Range("H18").Value = 1049.8572
$ws.Range("I18").Value = 712.25
$ws.Range("K18").Value = 2136.75
$ws.Range("M18").Value = -1967.75
# Row 113 (hunk 19)
$ws.Range("H113").Value = 1171.9166
$ws.Range("J113").Value = 1256.4
$ws.Range("L113").Value = 3769.2
$ws.Range("N113").Value = -8109.200000000001
# Row 121 (hunk 20)
$ws.Range("H121").Value = 771.75
$ws.Range("I121").Value = 566
$ws.Range("K121").Value = 1698
$ws.Range("M121").Value = -388
# Row 131 (hunk 21)
$ws.Range("H131").Value = 5881.709
$ws.Range("I131").Value = 8884.833000000001
$ws.Range("J131").Value = 5513.9795
$ws.Range("K131").Value = 26654.499
$ws.Range("L131").Value = 16541.9385
$ws.Range("M131").Value = -21614.499
$ws.Range("N131").Value = -26621.9385
# Row 135 (hunk 22)
$ws.Range("H135").Value = 1833.6316
$ws.Range("I135").Value = 1889.9375
$ws.Range("K135").Value = 17009.4375
$ws.Range("M135").Value = -14474.4375
# Row 136 (hunk 23)
$ws.Range("H136").Value = 9924.333000000001
$ws.Range("I136").Value = 9924.333000000001
$ws.Range("K136").Value = 29772.999
$ws.Range("M136").Value = -24672.999

$ws = $wb.Worksheets.Item("GSM")
# Row 14 (hunk 24)
$ws.Range("H14").Value = 5010750
$ws.Range("J14").Value = 6680000
$ws.Range("L14").Value = 6680000
$ws.Range("N14").Value = -6680336
# Row 63 (hunk 25)
$ws.Range("H63").Value = 27200
$ws.Range("J63").Value = 27200
$ws.Range("L63").Value = 27200
$ws.Range("N63").Value = -28572
# Row 66 (hunk 26)
$ws.Range("H66").Value = 27200
$ws.Range("J66").Value = 27200
$ws.Range("L66").Value = 81600
$ws.Range("N66").Value = -88464
# Row 101 (hunk 27)
$ws.Range("H101").Value = 49999
$ws.Range("J101").Value = 49999
$ws.Range("L101").Value = 49999
$ws.Range("N101").Value = -56489
# Row 102 (hunk 28)
$ws.Range("H102").Value = 2066.3667
$ws.Range("I102").Value = 1405.3478
$ws.Range("J102").Value = 4238.2856
$ws.Range("K102").Value = 1405.3478
$ws.Range("L102").Value = 4238.2856
$ws.Range("M102").Value = 216.6522
$ws.Range("N102").Value = -7482.2856

$ws = $wb.Worksheets.Item("LTW")
# Row 11 (hunk 29)
$ws.Range("H11").Value = 81.666664
$ws.Range("I11").Value = 110
$ws.Range("K11").Value = 110
$ws.Range("M11").Value = 30
# Row 55 (hunk 30)
$ws.Range("H55").Value = 22727750
$ws.Range("I55").Value = 28571870
$ws.Range("K55").Value = 28571870
$ws.Range("M55").Value = -28571697
# Row 61 (hunk 31)
$ws.Range("H61").Value = 1884.5555
$ws.Range("I61").Value = 1838.2
$ws.Range("J61").Value = 1942.5
$ws.Range("K61").Value = 1838.2
$ws.Range("L61").Value = 1942.5
$ws.Range("M61").Value = -1636.2
$ws.Range("N61").Value = -2346.5
# Row 64 (hunk 32)
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = ""
# Row 67 (hunk 33)
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = ""
# Row 100 (hunk 34)
$ws.Range("H100").Value = 3124
$ws.Range("I100").Value = 2264
$ws.Range("K100").Value = 2264
$ws.Range("M100").Value = -1723
# Row 108 (hunk 35)
$ws.Range("H108").Value = 73995
$ws.Range("J108").Value = 73995
$ws.Range("L108").Value = 73995
$ws.Range("N108").Value = -81675
# Row 109 (hunk 36)
$ws.Range("H109").Value = 98016
$ws.Range("J109").Value = 98016
$ws.Range("L109").Value = 98016
$ws.Range("N109").Value = -100790
# Row 113 (hunk 37)
$ws.Range("H113").Value = 1884.5555
$ws.Range("I113").Value = 1838.2
$ws.Range("J113").Value = 1942.5
$ws.Range("K113").Value = 1838.2
$ws.Range("L113").Value = 1942.5
$ws.Range("M113").Value = 331.8
$ws.Range("N113").Value = -6282.5
# Row 123 (hunk 38)
$ws.Range("H123").Value = 87995
$ws.Range("J123").Value = 87995
$ws.Range("L123").Value = 87995
$ws.Range("N123").Value = -97795
# Row 136 (hunk 39)
$ws.Range("H136").Value = 100008.734
$ws.Range("I136").Value = 18482.334
$ws.Range("K136").Value = 55447.00199999999
$ws.Range("M136").Value = -52897.00199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 15 (hunk 40)
$ws.Range("H15").Value = 80007
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = ""
# Row 63 (hunk 41)
$ws.Range("H63").Value = 34749.6
$ws.Range("J63").Value = 49999.332
$ws.Range("L63").Value = 49999.332
$ws.Range("N63").Value = -51247.332
# Row 66 (hunk 42)
$ws.Range("H66").Value = 34749.6
$ws.Range("J66").Value = 49999.332
$ws.Range("L66").Value = 149997.996
$ws.Range("N66").Value = -156237.996
# Row 96 (hunk 43)
$ws.Range("H96").Value = 4999
$ws.Range("I96").Value = 4999
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 4999
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -3626
$ws.Range("N96").Value = ""
# Row 100 (hunk 44)
$ws.Range("H100").Value = 1318.3125
$ws.Range("I100").Value = 3546.3333
$ws.Range("J100").Value = 804.1539
$ws.Range("K100").Value = 7092.6666
$ws.Range("L100").Value = 1608.3078
$ws.Range("M100").Value = -6551.6666
$ws.Range("N100").Value = -2690.3078
# Row 103 (hunk 45)
$ws.Range("H103").Value = 60555.5
$ws.Range("J103").Value = 60555.5
$ws.Range("L103").Value = 60555.5
$ws.Range("N103").Value = -62899.5
# Row 112 (hunk 46)
$ws.Range("H112").Value = 101255.29
$ws.Range("J112").Value = 101255.29
$ws.Range("L112").Value = 101255.29
$ws.Range("N112").Value = -104209.29
# Row 123 (hunk 47)
$ws.Range("H123").Value = 74015.8
$ws.Range("J123").Value = 74015.8
$ws.Range("L123").Value = 74015.8
$ws.Range("N123").Value = -83815.8
# Row 125 (hunk 48)
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""
